# Update "江西-漫展信息.xlsx": the first upcoming event ("江西·高安首届静卿国风动漫文化展览会",
# previously row 2) was removed from the feed. All subsequent events shift up by one row
# (columns B..I), the running index in column A stays fixed (1..15), the last row (old row 17)
# is dropped, and a handful of ticket-count / price cells were refreshed with newer values.
# This affects both the "展览" sheet and the "全部类型" sheet (sheet index 1 and 4), which are
# mirror copies of each other. "演出" and "本地生活" (index 2 and 3) are untouched.

$wb = $excel.ActiveWorkbook

$sheetIndexes = @(1, 4)

foreach ($idx in $sheetIndexes) {
    $ws = $wb.Worksheets.Item($idx)

    # Shift the event data (columns B..I) up by one row: row 3..17 -> row 2..16.
    # Column A (the running index 1..15) is intentionally left untouched.
    $ws.Range("B3:I17").Copy($ws.Range("B2:I16"))

    # Drop the now-duplicated last row (old row 17), shrinking the sheet to 16 rows.
    $ws.Rows.Item(17).Delete()

    # A few cells received fresher values beyond the plain shift-up (ticket interest /
    # lowest price updates scraped at a later time).
    $ws.Range("G3").Value2 = "不可售"
    $ws.Range("G4").Value2 = "已售罄"
    $ws.Range("F5").Value2 = 4603
    $ws.Range("F7").Value2 = 385
    $ws.Range("F8").Value2 = 1361
    $ws.Range("F9").Value2 = 898
    $ws.Range("G9").Value2 = 55
    $ws.Range("F10").Value2 = 55
    $ws.Range("F11").Value2 = 1022
    $ws.Range("F13").Value2 = 562
    $ws.Range("F15").Value2 = 259
    $ws.Range("F16").Value2 = 23
}
